$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "teacher1"
$ws.Range("B6").Value = 1234
$ws.Range("C6").Value = "teacher"
$ws.Range("D6").Value = "CS-A"
$ws.Range("E6").Value = "math"

$ws.Range("A7").Value = "teacher1"
$ws.Range("B7").Value = 1234
$ws.Range("C7").Value = "teacher"
$ws.Range("D7").Value = "CS-B"
$ws.Range("E7").Value = "math"

$ws.Range("K10").Select()
